# Remove the "Knowledge Base" guidance paragraphs ("Please note - any
# references ... should be checked." / "This means that GPT has either
# extrapolated ... has to be checked manually.") from every slide's notes
# page. The notes body placeholder is left in place, just emptied out.

$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    $notes = $s.NotesPage

    # The notes body text lives in placeholder index 2 ("Notes Placeholder
    # 2") of the notes page - index 1 is the slide image placeholder.
    $notesBody = $null
    foreach ($shp in $notes.Shapes.Placeholders) {
        if ($shp.Name -like "Notes Placeholder*") {
            $notesBody = $shp
            break
        }
    }
    if ($notesBody -eq $null) {
        $notesBody = $notes.Shapes.Placeholders.Item(2)
    }

    $notesBody.TextFrame.TextRange.Text = ""
}
